# Rename the two worksheets and make the first sheet ("prime") the
# selected/active tab instead of the second one.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Name = "prime"
$ws2.Name = "emails"

# Make "prime" (sheet 1) the active/selected tab; this clears
# tabSelected on "emails" and moves it to "prime", and updates the
# workbook's remembered active tab back to the first sheet.
$ws1.Select()
$ws1.Activate()
